# Auto-generated edit script: update betting odds values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("J5").Value = 1.05
$ws.Range("K5").Value = 11
$ws.Range("N5").Value = 1.88
$ws.Range("O5").Value = 1.93

# Row 6
$ws.Range("H6").Value = 9
$ws.Range("I6").Value = 35
$ws.Range("T6").Value = 14
$ws.Range("U6").Value = 7.8
$ws.Range("W6").Value = 5.9
$ws.Range("AA6").Value = 30
$ws.Range("AC6").Value = 200
$ws.Range("AE6").Value = 200
$ws.Range("AF6").Value = 450
$ws.Range("AG6").Value = 200
$ws.Range("AJ6").Value = 500

# Row 7
$ws.Range("G7").Value = 3.05
$ws.Range("H7").Value = 3.35
$ws.Range("I7").Value = 2.07
$ws.Range("O7").Value = 1.83
$ws.Range("T7").Value = 8.5
$ws.Range("U7").Value = 13.5
$ws.Range("V7").Value = 9.25
$ws.Range("W7").Value = 30
$ws.Range("X7").Value = 21
$ws.Range("AA7").Value = 5.8
$ws.Range("AB7").Value = 11.25
$ws.Range("AC7").Value = 45
$ws.Range("AF7").Value = 8.75
$ws.Range("AG7").Value = 7.4
$ws.Range("AH7").Value = 15.5
$ws.Range("AI7").Value = 13
$ws.Range("AJ7").Value = 20

# Row 11
$ws.Range("L11").Value = 1.22
$ws.Range("M11").Value = 4

# Row 18
$ws.Range("H18").Value = 3.15
$ws.Range("N18").Value = 2
$ws.Range("O18").Value = 1.65
$ws.Range("P18").Value = 1.39
$ws.Range("Q18").Value = 2.45
$ws.Range("T18").Value = 5.7
$ws.Range("Z18").Value = 8.25
$ws.Range("AA18").Value = 5.4
$ws.Range("AC18").Value = 55
$ws.Range("AD18").Value = 400
$ws.Range("AE18").Value = 7.9

# Row 20
$ws.Range("L20").Value = 1.13
$ws.Range("M20").Value = 6
$ws.Range("N20").Value = 1.44
$ws.Range("O20").Value = 2.7

# Row 21
$ws.Range("H21").Value = 3.5
$ws.Range("J21").Value = 1.05
$ws.Range("K21").Value = 11
$ws.Range("L21").Value = 1.25
$ws.Range("M21").Value = 3.75
$ws.Range("N21").Value = 1.85
$ws.Range("O21").Value = 1.95
$ws.Range("R21").Value = 1.67
$ws.Range("S21").Value = 2.1
$ws.Range("T21").Value = 9
$ws.Range("U21").Value = 11
$ws.Range("Z21").Value = 12
$ws.Range("AE21").Value = 11

# Row 22
$ws.Range("G22").Value = 2.15
$ws.Range("H22").Value = 3.6
$ws.Range("I22").Value = 3.1
$ws.Range("J22").Value = 1.04
$ws.Range("K22").Value = 12
$ws.Range("N22").Value = 1.83
$ws.Range("O22").Value = 1.98
$ws.Range("W22").Value = 19
$ws.Range("AD22").Value = 201
$ws.Range("AF22").Value = 17

# Row 23
$ws.Range("G23").Value = 1.44
$ws.Range("I23").Value = 6
$ws.Range("AC23").Value = 51
$ws.Range("AG23").Value = 19

# Row 26
$ws.Range("G26").Value = 1.67
$ws.Range("I26").Value = 4.33
$ws.Range("Z26").Value = 17

# Row 29
$ws.Range("G29").Value = 2.55
$ws.Range("I29").Value = 2.88
$ws.Range("K29").Value = 8
$ws.Range("L29").Value = 1.4
$ws.Range("M29").Value = 2.75
$ws.Range("N29").Value = 2.25
$ws.Range("O29").Value = 1.62
$ws.Range("U29").Value = 11
$ws.Range("V29").Value = 10
$ws.Range("Y29").Value = 34

# Row 33
$ws.Range("G33").Value = 4.15
$ws.Range("I33").Value = 1.65
$ws.Range("K33").Value = 8.75
$ws.Range("L33").Value = 1.2
$ws.Range("M33").Value = 4.05
$ws.Range("N33").Value = 1.6
$ws.Range("O33").Value = 2.18
$ws.Range("P33").Value = 1.32
$ws.Range("Q33").Value = 3.15
$ws.Range("R33").Value = 1.65
$ws.Range("S33").Value = 2.1
$ws.Range("T33").Value = 15
$ws.Range("U33").Value = 26
$ws.Range("V33").Value = 14
$ws.Range("W33").Value = 65
$ws.Range("X33").Value = 35
$ws.Range("Y33").Value = 37
$ws.Range("Z33").Value = 8.75
$ws.Range("AA33").Value = 8
$ws.Range("AB33").Value = 14.5
$ws.Range("AC33").Value = 55
$ws.Range("AD33").Value = 350
$ws.Range("AE33").Value = 8.75
$ws.Range("AF33").Value = 8.75
$ws.Range("AH33").Value = 13
$ws.Range("AJ33").Value = 21

# Row 34
$ws.Range("G34").Value = 2.82
$ws.Range("H34").Value = 3.25
$ws.Range("I34").Value = 2.27
$ws.Range("J34").Value = 1.06
$ws.Range("K34").Value = 7.3
$ws.Range("L34").Value = 1.3
$ws.Range("M34").Value = 3.2
$ws.Range("N34").Value = 1.9
$ws.Range("O34").Value = 1.82
$ws.Range("P34").Value = 1.42
$ws.Range("Q34").Value = 2.67
$ws.Range("R34").Value = 1.72
$ws.Range("S34").Value = 2
$ws.Range("T34").Value = 9
$ws.Range("U34").Value = 14.5
$ws.Range("V34").Value = 10.5
$ws.Range("W34").Value = 35
$ws.Range("X34").Value = 25
$ws.Range("Y34").Value = 32
$ws.Range("Z34").Value = 7.3
$ws.Range("AA34").Value = 6.4
$ws.Range("AB34").Value = 13.5
$ws.Range("AC34").Value = 60
$ws.Range("AD34").Value = 450
$ws.Range("AE34").Value = 8.25
$ws.Range("AF34").Value = 11.5
$ws.Range("AG34").Value = 9
$ws.Range("AH34").Value = 23
$ws.Range("AI34").Value = 18
$ws.Range("AJ34").Value = 27
